# Implementing add_assessment_result_to_matrix and some other minor changes
#
# The "Use case template" sheet tracked three baseline use cases
# (Baseline-SDMFFP1/2/3) in columns B:D. This drops the third baseline
# use case column (D) entirely, renames the remaining two headers from
# "Baseline-SDMFFPn" to "SDMFFPn", and records the assertion results that
# now apply to each of the two remaining use cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Use case template")

# Drop the third baseline column (Baseline-SDMFFP3) - shifts dimension
# from A1:D44 down to A1:C44 and removes the now-unused shared string.
$ws.Columns.Item(4).Delete() | Out-Null

# Rename the remaining two use-case headers.
$ws.Range("B1").Value = "SDMFFP1"
$ws.Range("C1").Value = "SDMFFP2"

# Record which assertion results feed into each use case's matrix.
$ws.Range("B7").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("B21").Value = 1
$ws.Range("C22").Value = 1
$ws.Range("C24").Value = 1

$ws.Range("D24").Select() | Out-Null
